$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 200, shifting existing rows 200:276 down to 202:278
$ws.Rows("200:201").Insert()

# New row 200 data
$ws.Range("A200").Value = 7
$ws.Range("B200").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C200").Value = "Ñuble"
$ws.Range("D200").Value = 44704
$ws.Range("E200").Value = 16
$ws.Range("F200").Value = 100112008
$ws.Range("G200").Value = "Coliflor"
$ws.Range("H200").Value = "Sin especificar"
$ws.Range("I200").Value = "Primera"
$ws.Range("J200").Value = 200
$ws.Range("K200").Value = 900
$ws.Range("L200").Value = 1000
$ws.Range("M200").Value = 950
$ws.Range("N200").Value = "$/unidad"
$ws.Range("O200").Value = "Provincia de Diguillín"
$ws.Range("P200").Value = 950
$ws.Range("Q200").Value = 1
$ws.Range("R200").Value = "Hortaliza"

# New row 201 data
$ws.Range("A201").Value = 7
$ws.Range("B201").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C201").Value = "Ñuble"
$ws.Range("D201").Value = 44704
$ws.Range("E201").Value = 16
$ws.Range("F201").Value = 100112008
$ws.Range("G201").Value = "Coliflor"
$ws.Range("H201").Value = "Sin especificar"
$ws.Range("I201").Value = "Segunda"
$ws.Range("J201").Value = 60
$ws.Range("K201").Value = 800
$ws.Range("L201").Value = 800
$ws.Range("M201").Value = 800
$ws.Range("N201").Value = "$/unidad"
$ws.Range("O201").Value = "Provincia de Diguillín"
$ws.Range("P201").Value = 800
$ws.Range("Q201").Value = 1
$ws.Range("R201").Value = "Hortaliza"

# Ensure D column on new rows keeps the date-style number format used elsewhere in column D
$ws.Range("D200:D201").NumberFormat = $ws.Range("D199").NumberFormat
